# Remove the obsolete {gender} placeholder (and the trailing space that
# separated it from "гражданин") from the "Заказчик"/"Исполнитель" clause.
$d = $word.ActiveDocument

$d.Content.Find.Execute("{gender} ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2)
